# Add `0` to vaccination dose counts for "12-17" (C) and "11-" (B) age
# groups on both the "1st dose" and "2nd dose" sheets, matching rows
# that already had data for the other age groups.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "1st dose"
$ws2 = $wb.Worksheets.Item(2)   # "2nd dose"

# --- Sheet 1: "1st dose" ---------------------------------------------

# Header cells already carry style s="4"; set their text to the
# existing shared-string values ("11-" for B, "12-17" for C).
$ws1.Range("B1").Value = "11-"
$ws1.Range("C1").Value = "12-17"

# Rows 2-21: B and C cells already exist (empty, style s="2") -> fill
# with 0.
$ws1.Range("B2:B21").Value = 0
$ws1.Range("C2:C17").Value = 0

# Rows 22-26: B cells don't exist yet on this sheet; copy the number
# format from the neighboring C cell (same row/style) before setting
# the value so the new cell gets style s="2" like the rest of the
# column.
$ws1.Range("B22:B26").NumberFormat = $ws1.Range("C22").NumberFormat
$ws1.Range("B22:B26").Value = 0

# --- Sheet 2: "2nd dose" ----------------------------------------------

# Rows 11-21: B cells already exist (empty, style s="2") -> fill with 0.
$ws2.Range("B11:B21").Value = 0
$ws2.Range("C11:C17").Value = 0

# Rows 22-26: B cells don't exist yet on this sheet either; copy the
# number format from the neighboring C cell first so the new cell
# picks up style s="2".
$ws2.Range("B22:B26").NumberFormat = $ws2.Range("C22").NumberFormat
$ws2.Range("B22:B26").Value = 0

# --- View state ---------------------------------------------------------
# Sheet2 ends up showing the B11:B26 selection (not the active tab).
$ws2.Activate()
$ws2.Range("B11:B26").Select()

# Sheet1 becomes the active tab, scrolled back to the top with C19
# selected.
$ws1.Activate()
$ws1.Range("C19").Select()
